$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ticker")

# Capture a known "default" (unstyled) cell style to restore onto cells after
# forcing a text (quote-prefixed) entry below - keeps the output style-free,
# matching the original workbook's formatting.
$defaultStyle = $ws.Cells.Item(1, 1).Style

function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $defaultStyle
}

# Remove the "link" column (column E) entirely - shifts dimension from A1:E10 to A1:D10
$ws.Columns.Item(5).Delete()

# Update Gold row (row 4): percentage/number change now flat ("UNCH"), quote updated
Set-TextValue $ws.Cells.Item(4, 1) "0"
$ws.Cells.Item(4, 2).Value = "UNCH "
Set-TextValue $ws.Cells.Item(4, 3) "1,211.20"

# Update Oil row (row 5): percentage/number change now flat ("UNCH"), quote updated
Set-TextValue $ws.Cells.Item(5, 1) "0"
$ws.Cells.Item(5, 2).Value = "UNCH "
Set-TextValue $ws.Cells.Item(5, 3) "48.03"

# Update EUR/USD row (row 6): change flipped to positive, quote updated
Set-TextValue $ws.Cells.Item(6, 1) "+0.07%"
Set-TextValue $ws.Cells.Item(6, 2) "0.0007 "
Set-TextValue $ws.Cells.Item(6, 3) "1.0631"
